# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22),
# keeping the existing date formatting/style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45221
}
